$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'mens compression basketball pants'
$ws.Cells.Item(2, 1).Value = 'mens basketball tights with knee pads'
$ws.Cells.Item(3, 1).Value = 'capri pants for boys'
$ws.Cells.Item(4, 1).Value = 'knee pads with leggings'
$ws.Cells.Item(5, 1).Value = 'compression pants for men'
$ws.Cells.Item(6, 1).Value = 'compression mens leggings'
$ws.Cells.Item(7, 1).Value = 'basketball protective gear'
$ws.Cells.Item(8, 1).Value = 'basketball spandex pants'
$ws.Cells.Item(9, 1).Value = 'compression tights girls'
$ws.Cells.Item(10, 1).Value = 'wrestling pants'
$ws.Cells.Item(11, 1).Value = 'mens compression capris'
$ws.Cells.Item(12, 1).Value = 'baseball sliding tights'
$ws.Cells.Item(13, 1).Value = 'boys basketball compression pants with knee pads'
$ws.Cells.Item(14, 1).Value = 'compression men capri'
$ws.Cells.Item(15, 1).Value = 'volleyball knee pads men'
$ws.Cells.Item(16, 1).Value = 'mens basketball'
$ws.Cells.Item(17, 1).Value = 'athletic leggings boys'
$ws.Cells.Item(18, 1).Value = 'knee pad weightlifting'
$ws.Cells.Item(19, 1).Value = 'basketball guide'
$ws.Cells.Item(20, 1).Value = 'basketball knee pads boys'
$ws.Cells.Item(21, 1).Value = 'youth padded compression pants basketball'
$ws.Cells.Item(22, 1).Value = 'leggings for men'
$ws.Cells.Item(23, 1).Value = 'cycling capris'
$ws.Cells.Item(24, 1).Value = 'padded knee pads for basketball'
$ws.Cells.Item(25, 1).Value = 'pant with knee pad'
$ws.Cells.Item(26, 1).Value = 'tight pants for men'
$ws.Cells.Item(27, 1).Value = 'basketball kneepads'
$ws.Cells.Item(28, 1).Value = 'thigh compression leggings'
$ws.Cells.Item(29, 1).Value = 'compression pants men pack'
$ws.Cells.Item(30, 1).Value = 'hiking knee pads'
$ws.Cells.Item(31, 1).Value = 'sliding pants'
$ws.Cells.Item(32, 1).Value = 'boys compression pants'
$ws.Cells.Item(33, 1).Value = 'cycling knee pads for men'
$ws.Cells.Item(34, 1).Value = 'youth boys knee pads for basketball'
$ws.Cells.Item(35, 1).Value = 'mens spandex leggings'
$ws.Cells.Item(36, 1).Value = 'compression knee padded'
$ws.Cells.Item(37, 1).Value = 'youth basketball kneepads'
$ws.Cells.Item(38, 1).Value = 'baseball equipment for boys'
$ws.Cells.Item(39, 1).Value = 'youth football tights boys'
$ws.Cells.Item(40, 1).Value = 'men compression workout pants'
$ws.Cells.Item(41, 1).Value = 'padded knee basketball'
$ws.Cells.Item(42, 1).Value = 'mens compression pants running'
$ws.Cells.Item(43, 1).Value = 'basketball youth leggings'
$ws.Cells.Item(44, 1).Value = 'compression basketball leggings youth'
$ws.Cells.Item(45, 1).Value = 'men compression running pants'
$ws.Cells.Item(46, 1).Value = 'wrestling knee pads youth'
$ws.Cells.Item(47, 1).Value = 'mens leggings sports'
$ws.Cells.Item(48, 1).Value = 'skin leggings men'
$ws.Cells.Item(49, 1).Value = 'youth football pants'
$ws.Cells.Item(50, 1).Value = 'mens leggings for sports'
$ws.Cells.Item(51, 1).Value = 'black boys compression pants'
$ws.Cells.Item(52, 1).Value = 'basketball team clothes'
$ws.Cells.Item(53, 1).Value = 'men compression tights'
$ws.Cells.Item(54, 1).Value = 'youth spandex leggings boys'
$ws.Cells.Item(55, 1).Value = 'baseball gear for men'
$ws.Cells.Item(56, 1).Value = 'sports knee pads'
$ws.Cells.Item(57, 1).Value = 'leggings men'
$ws.Cells.Item(58, 1).Value = 'pain in hip down leg'
$ws.Cells.Item(59, 1).Value = 'compression men tights'
$ws.Cells.Item(60, 1).Value = 'softball pants for men'
$ws.Cells.Item(61, 1).Value = 'youth basketball knee pads for boys'
$ws.Cells.Item(62, 1).Value = 'mens running tights'
$ws.Cells.Item(63, 1).Value = 'youth baseball gear'
$ws.Cells.Item(64, 1).Value = 'football pants adult xl'
$ws.Cells.Item(65, 1).Value = 'padded knee compression'
$ws.Cells.Item(66, 1).Value = 'compression pants men 3/4'
$ws.Cells.Item(67, 1).Value = 'sliding pants softball youth'
$ws.Cells.Item(68, 1).Value = 'hockey hip pads'
$ws.Cells.Item(69, 1).Value = 'compression knees'
$ws.Cells.Item(70, 1).Value = 'compression men'
$ws.Cells.Item(71, 1).Value = 'compression mens tights'
$ws.Cells.Item(72, 1).Value = 'compression tights men'
$ws.Cells.Item(73, 1).Value = 'legging for men'
$ws.Cells.Item(74, 1).Value = 'football knee pads adult'
$ws.Cells.Item(75, 1).Value = 'football pants with pads mens'
$ws.Cells.Item(76, 1).Value = 'mens basketball tights and leggings'
$ws.Cells.Item(77, 1).Value = 'knee protector'
$ws.Cells.Item(78, 1).Value = 'basketball apparel mens'
$ws.Cells.Item(79, 1).Value = 'adult knee pads'
$ws.Cells.Item(80, 1).Value = 'black softball pants girls'
$ws.Cells.Item(81, 1).Value = 'youth football pads for pants'
$ws.Cells.Item(82, 1).Value = 'youth black compression pants'
$ws.Cells.Item(83, 1).Value = '3/4 compression pants men'
$ws.Cells.Item(84, 1).Value = 'girls hiking pants'
$ws.Cells.Item(85, 1).Value = 'youth tights boys basketball'
$ws.Cells.Item(86, 1).Value = 'down pants men'
$ws.Cells.Item(87, 1).Value = 'basketball volleyball knee pads'
$ws.Cells.Item(88, 1).Value = 'legging for men sport'
$ws.Cells.Item(89, 1).Value = 'polyester hex mesh'
$ws.Cells.Item(90, 1).Value = 'basketball knee pad tights'
$ws.Cells.Item(91, 1).Value = 'capris tights'
$ws.Cells.Item(92, 1).Value = 'thigh compression pants'
$ws.Cells.Item(93, 1).Value = 'bjj pants men'
$ws.Cells.Item(94, 1).Value = 'wrestling pads'
$ws.Cells.Item(95, 1).Value = 'basketball athletic tights'
$ws.Cells.Item(96, 1).Value = 'wrestling knee pads youth 2 pack'
$ws.Cells.Item(97, 1).Value = 'compression for knees'
$ws.Cells.Item(98, 1).Value = 'youth football pants with pads black'
$ws.Cells.Item(99, 1).Value = 'girls knee pads'
$ws.Cells.Item(100, 1).Value = 'compression pants football'
